# Auto-generated edit script applying the Cactuar_Profits.xlsx price-refresh diff.
# Each hunk corresponds to a single leve row whose market-price-derived columns
# (H, I, J, K, L, M, N) were refreshed by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1444.0212
$ws.Range("J129").Value = 2338.4614
$ws.Range("L129").Value = 7015.3842
$ws.Range("N129").Value = -17015.3842
$ws.Range("H138").Value = 2193.2
$ws.Range("J138").Value = 2271.1177
$ws.Range("L138").Value = 6813.353099999999
$ws.Range("N138").Value = -17093.3531
$ws.Range("H141").Value = 4744.1113
$ws.Range("I141").Value = 2874.75
$ws.Range("K141").Value = 8624.25
$ws.Range("M141").Value = -3444.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2202.0244
$ws.Range("I45").Value = 1814.0968
$ws.Range("J45").Value = 3404.6
$ws.Range("K45").Value = 1814.0968
$ws.Range("L45").Value = 3404.6
$ws.Range("M45").Value = -1437.0968
$ws.Range("N45").Value = -4158.6
$ws.Range("H61").Value = 5450.486
$ws.Range("I61").Value = 5031.7666
$ws.Range("J61").Value = 7962.8
$ws.Range("K61").Value = 5031.7666
$ws.Range("L61").Value = 7962.8
$ws.Range("M61").Value = -4819.7666
$ws.Range("N61").Value = -8386.799999999999
$ws.Range("H74").Value = 14707307
$ws.Range("I74").Value = 17858214
$ws.Range("J74").Value = 3069.6667
$ws.Range("K74").Value = 17858214
$ws.Range("L74").Value = 3069.6667
$ws.Range("M74").Value = -17857340
$ws.Range("N74").Value = -4817.6667
$ws.Range("H77").Value = 14707307
$ws.Range("I77").Value = 17858214
$ws.Range("J77").Value = 3069.6667
$ws.Range("K77").Value = 89291070
$ws.Range("L77").Value = 15348.3335
$ws.Range("M77").Value = -89286702
$ws.Range("N77").Value = -24084.3335
$ws.Range("H102").Value = 2349.0833
$ws.Range("I102").Value = 2308.0908
$ws.Range("K102").Value = 2308.0908
$ws.Range("M102").Value = -686.0907999999999
$ws.Range("H110").Value = 1461.9143
$ws.Range("I110").Value = 1239.8928
$ws.Range("J110").Value = 2350
$ws.Range("K110").Value = 1239.8928
$ws.Range("L110").Value = 2350
$ws.Range("M110").Value = 805.1071999999999
$ws.Range("N110").Value = -6440
$ws.Range("H122").Value = 1597.3611
$ws.Range("I122").Value = 1392.3438
$ws.Range("K122").Value = 4177.0314
$ws.Range("M122").Value = -1727.0314
$ws.Range("H132").Value = 11255.183
$ws.Range("I132").Value = 11147.933
$ws.Range("K132").Value = 33443.799
$ws.Range("M132").Value = -30913.799
$ws.Range("H136").Value = 5450.486
$ws.Range("I136").Value = 5031.7666
$ws.Range("J136").Value = 7962.8
$ws.Range("K136").Value = 15095.2998
$ws.Range("L136").Value = 23888.4
$ws.Range("M136").Value = -12545.2998
$ws.Range("N136").Value = -28988.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5032.452
$ws.Range("I20").Value = 6249.25
$ws.Range("J20").Value = 4283.654
$ws.Range("K20").Value = 6249.25
$ws.Range("L20").Value = 4283.654
$ws.Range("M20").Value = -6002.25
$ws.Range("N20").Value = -4777.654
$ws.Range("H37").Value = 2928.1667
$ws.Range("I37").Value = 513.8
$ws.Range("K37").Value = 513.8
$ws.Range("M37").Value = -376.8
$ws.Range("H69").Value = 119995
$ws.Range("J69").Value = 119995
$ws.Range("L69").Value = 119995
$ws.Range("N69").Value = -121617
$ws.Range("H72").Value = 119995
$ws.Range("J72").Value = 119995
$ws.Range("L72").Value = 359985
$ws.Range("N72").Value = -368097
$ws.Range("H75").Value = 29434.334
$ws.Range("I75").Value = 11999.5
$ws.Range("K75").Value = 11999.5
$ws.Range("M75").Value = -11063.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H78").Value = 29434.334
$ws.Range("I78").Value = 11999.5
$ws.Range("K78").Value = 35998.5
$ws.Range("M78").Value = -31318.5
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H86").Value = 3103.0605
$ws.Range("I86").Value = 3253.6
$ws.Range("J86").Value = 2977.611
$ws.Range("K86").Value = 3253.6
$ws.Range("L86").Value = 2977.611
$ws.Range("M86").Value = -2130.6
$ws.Range("N86").Value = -5223.611
$ws.Range("H88").Value = 18520
$ws.Range("J88").Value = 18520
$ws.Range("L88").Value = 18520
$ws.Range("N88").Value = -19332
$ws.Range("H89").Value = 3103.0605
$ws.Range("I89").Value = 3253.6
$ws.Range("J89").Value = 2977.611
$ws.Range("K89").Value = 16268
$ws.Range("L89").Value = 14888.055
$ws.Range("M89").Value = -10652
$ws.Range("N89").Value = -26120.055
$ws.Range("H91").Value = 18520
$ws.Range("J91").Value = 18520
$ws.Range("L91").Value = 18520
$ws.Range("N91").Value = -21328
$ws.Range("H105").Value = 2949.6086
$ws.Range("I105").Value = 2969.8
$ws.Range("J105").Value = 2815
$ws.Range("K105").Value = 2969.8
$ws.Range("L105").Value = 2815
$ws.Range("M105").Value = -1222.8
$ws.Range("N105").Value = -6309
$ws.Range("H107").Value = 2419.2632
$ws.Range("I107").Value = 2026.9143
$ws.Range("K107").Value = 2026.9143
$ws.Range("M107").Value = -106.9142999999999
$ws.Range("H134").Value = 3770.0476
$ws.Range("I134").Value = 3770.0476
$ws.Range("K134").Value = 11310.1428
$ws.Range("M134").Value = -8775.1428
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24393602
$ws.Range("I31").Value = 32260382
$ws.Range("J31").Value = 6584.3
$ws.Range("K31").Value = 32260382
$ws.Range("L31").Value = 6584.3
$ws.Range("M31").Value = -32260087
$ws.Range("N31").Value = -7174.3
$ws.Range("H34").Value = 24393602
$ws.Range("I34").Value = 32260382
$ws.Range("J34").Value = 6584.3
$ws.Range("K34").Value = 32260382
$ws.Range("L34").Value = 6584.3
$ws.Range("M34").Value = -32260180
$ws.Range("N34").Value = -6988.3
$ws.Range("H122").Value = 1773.6072
$ws.Range("I122").Value = 1569.2084
$ws.Range("K122").Value = 4707.6252
$ws.Range("M122").Value = -2257.6252
$ws.Range("H132").Value = 41669304
$ws.Range("I132").Value = 60608204
$ws.Range("K132").Value = 181824612
$ws.Range("M132").Value = -181822082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 38662104
$ws.Range("I4").Value = 35027016
$ws.Range("J4").Value = 49203856
$ws.Range("K4").Value = 105081048
$ws.Range("L4").Value = 147611568
$ws.Range("M4").Value = -105080936
$ws.Range("N4").Value = -147611792
$ws.Range("H107").Value = 899.8
$ws.Range("J107").Value = 999.75
$ws.Range("L107").Value = 2999.25
$ws.Range("N107").Value = -6839.25
$ws.Range("H129").Value = 2182.7827
$ws.Range("I129").Value = 961.0909
$ws.Range("J129").Value = 3302.6667
$ws.Range("K129").Value = 2883.2727
$ws.Range("L129").Value = 9908.000100000001
$ws.Range("M129").Value = 2116.7273
$ws.Range("N129").Value = -19908.0001
$ws.Range("H136").Value = 3040.6
$ws.Range("H141").Value = 5915.25
$ws.Range("I141").Value = 3870.3845
$ws.Range("K141").Value = 11611.1535
$ws.Range("M141").Value = -6431.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 166.91667
$ws.Range("I2").Value = 95.63636
$ws.Range("K2").Value = 95.63636
$ws.Range("M2").Value = 17.36364
$ws.Range("H102").Value = 10412976
$ws.Range("I102").Value = 11595780
$ws.Range("J102").Value = 4299.8
$ws.Range("K102").Value = 11595780
$ws.Range("L102").Value = 4299.8
$ws.Range("M102").Value = -11594158
$ws.Range("N102").Value = -7543.8
$ws.Range("H122").Value = 260468.69
$ws.Range("I122").Value = 456626.9
$ws.Range("J122").Value = 6616.8823
$ws.Range("K122").Value = 1369880.7
$ws.Range("L122").Value = 19850.6469
$ws.Range("M122").Value = -1367430.7
$ws.Range("N122").Value = -24750.6469
$ws.Range("H123").Value = 40959.383
$ws.Range("J123").Value = 40959.383
$ws.Range("L123").Value = 40959.383
$ws.Range("N123").Value = -45859.383
$ws.Range("H132").Value = 127787.31
$ws.Range("I132").Value = 168466.5
$ws.Range("K132").Value = 505399.5
$ws.Range("M132").Value = -502869.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1521.1034
$ws.Range("I16").Value = 917.2083
$ws.Range("J16").Value = 4419.8
$ws.Range("K16").Value = 917.2083
$ws.Range("L16").Value = 4419.8
$ws.Range("M16").Value = -747.2083
$ws.Range("N16").Value = -4759.8
$ws.Range("H93").Value = 1820.0476
$ws.Range("I93").Value = 1785.6666
$ws.Range("K93").Value = 1785.6666
$ws.Range("M93").Value = -537.6666
$ws.Range("H122").Value = 7849.9736
$ws.Range("I122").Value = 3844.2354
$ws.Range("K122").Value = 11532.7062
$ws.Range("M122").Value = -9082.706200000001
$ws.Range("H132").Value = 5489.7183
$ws.Range("I132").Value = 4795.5093
$ws.Range("K132").Value = 14386.5279
$ws.Range("M132").Value = -11856.5279

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 22000
$ws.Range("I21").Value = 22000
$ws.Range("K21").Value = 22000
$ws.Range("M21").Value = -21765
$ws.Range("H35").Value = 22000
$ws.Range("I35").Value = 22000
$ws.Range("K35").Value = 22000
$ws.Range("M35").Value = -21710
$ws.Range("H122").Value = 3789.5483
$ws.Range("J122").Value = 4544.7144
$ws.Range("L122").Value = 13634.1432
$ws.Range("N122").Value = -18534.1432
$ws.Range("H126").Value = 142858000
$ws.Range("I126").Value = 333333980
$ws.Range("K126").Value = 1000001940
$ws.Range("M126").Value = -999999470
$ws.Range("H132").Value = 7248368.5
$ws.Range("I132").Value = 12346778
$ws.Range("K132").Value = 37040334
$ws.Range("M132").Value = -37037804
$ws.Range("H136").Value = 5297.114
$ws.Range("I136").Value = 3591.4285
$ws.Range("K136").Value = 10774.2855
$ws.Range("M136").Value = -8224.2855
